$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Reset the "thick-bottom" data rows (3-11) to a uniform, non-wrapped height ---
for ($r = 3; $r -le 11; $r++) {
    $ws.Rows($r).RowHeight = 15.75
}

# --- Add the new note row 13 ---
# A13: reuse the existing header font (Arial 12, same as A3:A11) but without the
# heavy border / wrap-alignment those header cells carry, then paint it yellow.
$ws.Cells.Item(3, 1).Copy($ws.Cells.Item(13, 1))
$ws.Cells.Item(13, 1).Borders.LineStyle = -4142   # xlLineStyleNone
$ws.Cells.Item(13, 1).WrapText = $false
$ws.Cells.Item(13, 1).VerticalAlignment = -4107   # xlBottom (default/general)
$ws.Cells.Item(13, 1).Value = "***QCP has an error just before the hour is complete giving incorrect values"

# B13:G13 + A13 all get the same yellow fill
$ws.Range("A13:G13").Interior.Color = 65535
$ws.Rows(13).RowHeight = 15.75

# --- Selection / view state ---
$excel.Goto($ws.Range("A13:G13"), $true)
